# The course schema no longer tracks a "Dates" column.
# Remove the entire "Dates" column (column E: CourseCode, CourseName,
# Department, Instructor, Dates, Seat, Time, Venue, Description) and
# shift the remaining columns (Seat, Time, Venue, Description) left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").EntireColumn.Delete()
